# Applies the diff: duplicate "HIS (Routines)" into a new sheet named
# "HIS (Classes,Methods,Types)", adjust window/selection state, and make
# the new sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the first sheet (carries over all data/formatting) and place
# it right after the original.
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "HIS (Classes,Methods,Types)"

# Give each sheet the selection state seen in the target workbook.
$ws1.Activate() | Out-Null
$ws1.Range("A1:G38").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("J13").Select() | Out-Null

# The new sheet is the active / visible tab when the workbook is saved.
$ws2.Activate() | Out-Null
